$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4939.486
$ws.Range("I70").Value = 2147.9395
$ws.Range("J70").Value = 51000
$ws.Range("K70").Value = 6443.818499999999
$ws.Range("L70").Value = 153000
$ws.Range("M70").Value = -6173.818499999999
$ws.Range("N70").Value = -153540

$ws.Range("H73").Value = 4939.486
$ws.Range("I73").Value = 2147.9395
$ws.Range("J73").Value = 51000
$ws.Range("K73").Value = 6443.818499999999
$ws.Range("L73").Value = 153000
$ws.Range("M73").Value = -5507.818499999999
$ws.Range("N73").Value = -154872

$ws.Range("H88").Value = 71437020
$ws.Range("I88").Value = 166667400
$ws.Range("J88").Value = 14243.75
$ws.Range("K88").Value = 166667400
$ws.Range("L88").Value = 14243.75
$ws.Range("M88").Value = -166666994
$ws.Range("N88").Value = -15055.75

$ws.Range("H91").Value = 71437020
$ws.Range("I91").Value = 166667400
$ws.Range("J91").Value = 14243.75
$ws.Range("K91").Value = 166667400
$ws.Range("L91").Value = 14243.75
$ws.Range("M91").Value = -166665996
$ws.Range("N91").Value = -17051.75

$ws.Range("H100").Value = 2733.625
$ws.Range("I100").Value = 934.3333
$ws.Range("J100").Value = 3813.2
$ws.Range("K100").Value = 934.3333
$ws.Range("L100").Value = 3813.2
$ws.Range("M100").Value = -393.3333
$ws.Range("N100").Value = -4895.2

$ws.Range("H113").Value = 29037.818
$ws.Range("I113").Value = 61002
$ws.Range("K113").Value = 61002
$ws.Range("M113").Value = -57748

$ws.Range("H115").Value = 185
$ws.Range("I115").Value = 185
$ws.Range("K115").Value = 555
$ws.Range("M115").Value = 1012

$ws.Range("H137").Value = 2207.3462
$ws.Range("J137").Value = 2677.6155
$ws.Range("L137").Value = 8032.8465
$ws.Range("N137").Value = -13132.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1430.1428
$ws.Range("J45").Value = 1547.75
$ws.Range("L45").Value = 1547.75
$ws.Range("N45").Value = -2301.75

$ws.Range("H61").Value = 1803.1
$ws.Range("I61").Value = 1805.2858
$ws.Range("J61").Value = 1798
$ws.Range("K61").Value = 1805.2858
$ws.Range("L61").Value = 1798
$ws.Range("M61").Value = -1593.2858
$ws.Range("N61").Value = -2222

$ws.Range("H74").Value = 1123.0555
$ws.Range("I74").Value = 573.7931
$ws.Range("K74").Value = 573.7931
$ws.Range("M74").Value = 300.2069

$ws.Range("H77").Value = 1123.0555
$ws.Range("I77").Value = 573.7931
$ws.Range("K77").Value = 2868.9655
$ws.Range("M77").Value = 1499.0345

$ws.Range("H88").Value = 3395.7778
$ws.Range("I88").Value = 1663.75
$ws.Range("J88").Value = 4781.4
$ws.Range("K88").Value = 1663.75
$ws.Range("L88").Value = 4781.4
$ws.Range("M88").Value = -1257.75
$ws.Range("N88").Value = -5593.4

$ws.Range("H91").Value = 3395.7778
$ws.Range("I91").Value = 1663.75
$ws.Range("J91").Value = 4781.4
$ws.Range("K91").Value = 1663.75
$ws.Range("L91").Value = 4781.4
$ws.Range("M91").Value = -259.75
$ws.Range("N91").Value = -7589.4

$ws.Range("I97").Value = 966.6842
$ws.Range("J97").Value = 927.5
$ws.Range("K97").Value = 966.6842
$ws.Range("L97").Value = 927.5
$ws.Range("M97").Value = -470.6842
$ws.Range("N97").Value = -1919.5

$ws.Range("H136").Value = 1803.1
$ws.Range("I136").Value = 1805.2858
$ws.Range("J136").Value = 1798
$ws.Range("K136").Value = 5415.857400000001
$ws.Range("L136").Value = 5394
$ws.Range("M136").Value = -2865.857400000001
$ws.Range("N136").Value = -10494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1094250.5
$ws.Range("I86").Value = 1254232
$ws.Range("J86").Value = 667633.3
$ws.Range("K86").Value = 1254232
$ws.Range("L86").Value = 667633.3
$ws.Range("M86").Value = -1253109
$ws.Range("N86").Value = -669879.3

$ws.Range("H89").Value = 1094250.5
$ws.Range("I89").Value = 1254232
$ws.Range("J89").Value = 667633.3
$ws.Range("K89").Value = 6271160
$ws.Range("L89").Value = 3338166.5
$ws.Range("M89").Value = -6265544
$ws.Range("N89").Value = -3349398.5

$ws.Range("H99").Value = 1229
$ws.Range("I99").Value = 874.25
$ws.Range("K99").Value = 874.25
$ws.Range("M99").Value = 623.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3887.2856
$ws.Range("I58").Value = 2489.4
$ws.Range("J58").Value = 4663.8887
$ws.Range("K58").Value = 2489.4
$ws.Range("L58").Value = 4663.8887
$ws.Range("M58").Value = -2286.4
$ws.Range("N58").Value = -5069.8887

$ws.Range("H62").Value = 2273.375
$ws.Range("I62").Value = 2364.5
$ws.Range("K62").Value = 2364.5
$ws.Range("M62").Value = -1740.5

$ws.Range("H65").Value = 2273.375
$ws.Range("I65").Value = 2364.5
$ws.Range("K65").Value = 11822.5
$ws.Range("M65").Value = -8702.5

$ws.Range("H86").Value = 3733.2
$ws.Range("I86").Value = 3288.6667
$ws.Range("K86").Value = 3288.6667
$ws.Range("M86").Value = -2165.6667

$ws.Range("H89").Value = 3733.2
$ws.Range("I89").Value = 3288.6667
$ws.Range("K89").Value = 16443.3335
$ws.Range("M89").Value = -10827.3335

$ws.Range("H95").Value = 31000
$ws.Range("J95").Value = 31000
$ws.Range("L95").Value = 31000
$ws.Range("N95").Value = -36492

$ws.Range("H107").Value = 766.625
$ws.Range("I107").Value = 511
$ws.Range("K107").Value = 511
$ws.Range("M107").Value = 1409

$ws.Range("H134").Value = 2633.5417
$ws.Range("I134").Value = 2247.9524
$ws.Range("J134").Value = 5332.6665
$ws.Range("K134").Value = 6743.8572
$ws.Range("L134").Value = 15997.9995
$ws.Range("M134").Value = -4208.8572
$ws.Range("N134").Value = -21067.9995

$ws.Range("H136").Value = 3887.2856
$ws.Range("I136").Value = 2489.4
$ws.Range("J136").Value = 4663.8887
$ws.Range("K136").Value = 7468.200000000001
$ws.Range("L136").Value = 13991.6661
$ws.Range("M136").Value = -4918.200000000001
$ws.Range("N136").Value = -19091.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19590.648
$ws.Range("I131").Value = 720
$ws.Range("J131").Value = 21878
$ws.Range("K131").Value = 2160
$ws.Range("L131").Value = 65634
$ws.Range("M131").Value = 2880
$ws.Range("N131").Value = -75714

$ws.Range("H140").Value = 2567.6453
$ws.Range("I140").Value = 1314.4
$ws.Range("J140").Value = 3742.5625
$ws.Range("K140").Value = 3943.2
$ws.Range("L140").Value = 11227.6875
$ws.Range("M140").Value = 1236.8
$ws.Range("N140").Value = -21587.6875

$ws.Range("H141").Value = 3402.889
$ws.Range("I141").Value = 3660.8572
$ws.Range("J141").Value = 2500
$ws.Range("K141").Value = 10982.5716
$ws.Range("L141").Value = 7500
$ws.Range("M141").Value = -5802.571599999999
$ws.Range("N141").Value = -17860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 118.333336
$ws.Range("I107").Value = 122.5
$ws.Range("K107").Value = 122.5
$ws.Range("M107").Value = 1797.5

$ws.Range("H122").Value = 1095.7037
$ws.Range("I122").Value = 1111.84
$ws.Range("K122").Value = 3335.52
$ws.Range("M122").Value = -885.5199999999995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12832.923
$ws.Range("I40").Value = 13492.667
$ws.Range("K40").Value = 13492.667
$ws.Range("M40").Value = -13356.667

$ws.Range("H93").Value = 16667495
$ws.Range("I93").Value = 828.6429000000001
$ws.Range("J93").Value = 55556384
$ws.Range("K93").Value = 828.6429000000001
$ws.Range("L93").Value = 55556384
$ws.Range("M93").Value = 419.3570999999999
$ws.Range("N93").Value = -55558880

$ws.Range("H132").Value = 5892.636
$ws.Range("I132").Value = 5914.5713
$ws.Range("K132").Value = 17743.7139
$ws.Range("M132").Value = -15213.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2418.8
$ws.Range("I96").Value = 745
$ws.Range("J96").Value = 3534.6667
$ws.Range("K96").Value = 745
$ws.Range("L96").Value = 3534.6667
$ws.Range("M96").Value = 628
$ws.Range("N96").Value = -6280.6667
